$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 120, shifting the existing rows 120:212 down to 121:213
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new data point
$ws.Cells.Item(120, 1).Value  = 8
$ws.Cells.Item(120, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(120, 3).Value  = "Coquimbo"
$ws.Cells.Item(120, 4).Value  = 45280
$ws.Cells.Item(120, 5).Value  = 4
$ws.Cells.Item(120, 6).Value  = 100114007
$ws.Cells.Item(120, 7).Value  = "Jengibre"
$ws.Cells.Item(120, 8).Value  = "Sin especificar"
$ws.Cells.Item(120, 9).Value  = "Primera"
$ws.Cells.Item(120, 10).Value = 360
$ws.Cells.Item(120, 11).Value = 20000
$ws.Cells.Item(120, 12).Value = 22000
$ws.Cells.Item(120, 13).Value = 21000
$ws.Cells.Item(120, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(120, 15).Value = "Perú"
$ws.Cells.Item(120, 16).Value = 1615
$ws.Cells.Item(120, 17).Value = 13
$ws.Cells.Item(120, 18).Value = "Hortaliza"
